# Auto-generated Excel COM-interop script
# Applies numeric corrections to the per-job Leve profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 6: Days of Chunder / Antidote
$ws.Range("H6").Value = 3881.7273
$ws.Range("I6").Value = 283.16666
$ws.Range("J6").Value = 8200
$ws.Range("K6").Value = 849.4999799999999
$ws.Range("L6").Value = 24600
$ws.Range("M6").Value = -737.4999799999999
$ws.Range("N6").Value = -24824

# Row 33: Glazed and Confused / Clear Glass Lens
$ws.Range("H33").Value = 1561.8055
$ws.Range("I33").Value = 673.8823
$ws.Range("J33").Value = 2356.2632
$ws.Range("K33").Value = 673.8823
$ws.Range("L33").Value = 2356.2632
$ws.Range("M33").Value = -444.8823
$ws.Range("N33").Value = -2814.2632

# Row 113: Amaro Kart / Starch Glue
$ws.Range("H113").Value = 6877.1
$ws.Range("I113").Value = 5703.5
$ws.Range("J113").Value = 8637.5
$ws.Range("K113").Value = 5703.5
$ws.Range("L113").Value = 8637.5
$ws.Range("M113").Value = -2449.5
$ws.Range("N113").Value = -15145.5

# Row 138: All-night Crafting / Cunning Craftsman's Tisane
$ws.Range("H138").Value = 2509.2
$ws.Range("I138").Value = 2509.2
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 7527.599999999999
$ws.Range("L138").Value = 0
$ws.Range("M138").Value = -2387.599999999999
$ws.Range("N138").Value = ""

$ws = $wb.Worksheets.Item("ARM")
# Row 88: The Mast Chance / Adamantite Rivets
$ws.Range("H88").Value = 2425.077
$ws.Range("I88").Value = 1914.8572
$ws.Range("J88").Value = 3020.3333
$ws.Range("K88").Value = 1914.8572
$ws.Range("L88").Value = 3020.3333
$ws.Range("M88").Value = -1508.8572
$ws.Range("N88").Value = -3832.3333

# Row 91: The Rose and the Riveter (L) / Adamantite Rivets
$ws.Range("H91").Value = 2425.077
$ws.Range("I91").Value = 1914.8572
$ws.Range("J91").Value = 3020.3333
$ws.Range("K91").Value = 1914.8572
$ws.Range("L91").Value = 3020.3333
$ws.Range("M91").Value = -510.8571999999999
$ws.Range("N91").Value = -5828.3333

# Row 110: Scheduled Maintenance / Deepgold Ingot
$ws.Range("H110").Value = 4055.4443
$ws.Range("I110").Value = 3799.8572
$ws.Range("K110").Value = 3799.8572
$ws.Range("M110").Value = -1754.8572

# Row 114: A New Regular / Bluespirit Gauntlets of Fending
$ws.Range("H114").Value = 76000
$ws.Range("J114").Value = 76000
$ws.Range("L114").Value = 76000
$ws.Range("N114").Value = -84678

# Row 122: Haste for High Durium / High Durium Nugget
$ws.Range("H122").Value = 2498.1875
$ws.Range("I122").Value = 1909.8462
$ws.Range("K122").Value = 5729.5386
$ws.Range("M122").Value = -3279.5386

$ws = $wb.Worksheets.Item("BSM")
# Row 20: Smelt and Dealt / Iron Ingot
$ws.Range("H20").Value = 10482.588
$ws.Range("I20").Value = 12385.923
$ws.Range("J20").Value = 4296.75
$ws.Range("K20").Value = 12385.923
$ws.Range("L20").Value = 4296.75
$ws.Range("M20").Value = -12138.923
$ws.Range("N20").Value = -4790.75

# Row 105: Ingot to Wing It / Molybdenum Ingot
$ws.Range("H105").Value = 3723
$ws.Range("I105").Value = 4311.7856
$ws.Range("J105").Value = 2973.6365
$ws.Range("K105").Value = 4311.7856
$ws.Range("L105").Value = 2973.6365
$ws.Range("M105").Value = -2564.7856
$ws.Range("N105").Value = -6467.636500000001

# Row 107: The Gold Experience / Deepgold Nugget
$ws.Range("H107").Value = 1928.2222
$ws.Range("I107").Value = 1480.4615
$ws.Range("K107").Value = 1480.4615
$ws.Range("M107").Value = 439.5385000000001

$ws = $wb.Worksheets.Item("CRP")
# Row 16: Raise the Roof / Ash Lumber
$ws.Range("H16").Value = 921.75
$ws.Range("I16").Value = 733
$ws.Range("K16").Value = 733
$ws.Range("M16").Value = -446

# Row 17: Say It with Spears / Feathered Harpoon
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = ""
$ws.Range("N17").Value = ""

# Row 105: Zelkova, My Love / Zelkova Lumber
$ws.Range("H105").Value = 1636.5714
$ws.Range("I105").Value = 1576.25
$ws.Range("K105").Value = 1576.25
$ws.Range("M105").Value = 170.75

# Row 113: Patient Patients / White Ash Lumber
$ws.Range("H113").Value = 921.75
$ws.Range("I113").Value = 733
$ws.Range("K113").Value = 733
$ws.Range("M113").Value = 1437

# Row 132: Hull Lotta Damage / Ginseng Lumber
$ws.Range("H132").Value = 37932.36
$ws.Range("I132").Value = 31949.234
$ws.Range("K132").Value = 95847.702
$ws.Range("M132").Value = -93317.702

$ws = $wb.Worksheets.Item("CUL")
# Row 5: What a Sap / Maple Syrup
$ws.Range("H5").Value = 1027.0869
$ws.Range("J5").Value = 2990
$ws.Range("L5").Value = 8970
$ws.Range("N5").Value = -9194

# Row 121: A Cookie for Your Troubles / Coffee Biscuit
$ws.Range("H121").Value = 2842.5715
$ws.Range("I121").Value = 2546.3333
$ws.Range("J121").Value = 3375.8
$ws.Range("K121").Value = 7638.999899999999
$ws.Range("L121").Value = 10127.4
$ws.Range("M121").Value = -6328.999899999999
$ws.Range("N121").Value = -12747.4

# Row 135: Not-so-secret Ingredient / Royal Maple Syrup
$ws.Range("H135").Value = 1027.0869
$ws.Range("J135").Value = 2990
$ws.Range("L135").Value = 26910
$ws.Range("N135").Value = -31980

$ws = $wb.Worksheets.Item("GSM")
# Row 97: If I'd a Koppranickel for Every Time... / Koppranickel Ingot
$ws.Range("H97").Value = 968.4857
$ws.Range("J97").Value = 1346.1818
$ws.Range("L97").Value = 1346.1818
$ws.Range("N97").Value = -2338.1818

# Row 98: Cutting Deals / Durium Smallsword
$ws.Range("H98").Value = 29083.666
$ws.Range("J98").Value = 29083.666
$ws.Range("L98").Value = 29083.666
$ws.Range("N98").Value = -35073.666

# Row 107: Whetstones for the Workers / Hard Mudstone Whetstone
$ws.Range("H107").Value = 706.2222
$ws.Range("I107").Value = 618.0909
$ws.Range("J107").Value = 844.7143
$ws.Range("K107").Value = 618.0909
$ws.Range("L107").Value = 844.7143
$ws.Range("M107").Value = 1301.9091
$ws.Range("N107").Value = -4684.7143

# Row 113: Copious Crystal Cannons / Manasilver Nugget
$ws.Range("H113").Value = 224383.11
$ws.Range("J113").Value = 2313
$ws.Range("L113").Value = 2313
$ws.Range("N113").Value = -6653

$ws = $wb.Worksheets.Item("LTW")
# Row 7: Tan Before the Ban / Leather
$ws.Range("H7").Value = 12862.818
$ws.Range("I7").Value = 12862.818
$ws.Range("K7").Value = 12862.818
$ws.Range("M7").Value = -12750.818

# Row 22: Skin off Their Backs / Aldgoat Leather
$ws.Range("H22").Value = 3080
$ws.Range("J22").Value = 2750
$ws.Range("L22").Value = 2750
$ws.Range("N22").Value = -3340

# Row 27: Fire and Hide / Aldgoat Leather
$ws.Range("H27").Value = 3080
$ws.Range("J27").Value = 2750
$ws.Range("L27").Value = 2750
$ws.Range("N27").Value = -2964

# Row 61: Spelling Me Softly / Raptor Leather
$ws.Range("H61").Value = 1062.1666
$ws.Range("I61").Value = 1097.5714
$ws.Range("K61").Value = 1097.5714
$ws.Range("M61").Value = -895.5714

# Row 68: You Could Say It's a Moving Target / Wyvern Leather
$ws.Range("H68").Value = 2590.24
$ws.Range("I68").Value = 2259.432
$ws.Range("J68").Value = 5016.1665
$ws.Range("K68").Value = 2259.432
$ws.Range("L68").Value = 5016.1665
$ws.Range("M68").Value = -1510.432
$ws.Range("N68").Value = -6514.1665

# Row 71: They Call It Bloody Mary (L) / Wyvern Leather
$ws.Range("H71").Value = 2590.24
$ws.Range("I71").Value = 2259.432
$ws.Range("J71").Value = 5016.1665
$ws.Range("K71").Value = 11297.16
$ws.Range("L71").Value = 25080.8325
$ws.Range("M71").Value = -7553.16
$ws.Range("N71").Value = -32568.8325

# Row 82: Trainin' the Neck / Dragon Leather
$ws.Range("H82").Value = 1786.1875
$ws.Range("I82").Value = 832
$ws.Range("J82").Value = 3013
$ws.Range("K82").Value = 832
$ws.Range("L82").Value = 3013
$ws.Range("M82").Value = -471
$ws.Range("N82").Value = -3735

# Row 85: Training Is Only Skintight (L) / Dragon Leather
$ws.Range("H85").Value = 1786.1875
$ws.Range("I85").Value = 832
$ws.Range("J85").Value = 3013
$ws.Range("K85").Value = 832
$ws.Range("L85").Value = 3013
$ws.Range("M85").Value = 416
$ws.Range("N85").Value = -5509

# Row 110: Breeches of Trust / Gliderskin Breeches of Fending
$ws.Range("H110").Value = 87500
$ws.Range("J110").Value = 87500
$ws.Range("L110").Value = 87500
$ws.Range("N110").Value = -95680

# Row 113: Peace in Rest / Atrociraptor Leather
$ws.Range("H113").Value = 1062.1666
$ws.Range("I113").Value = 1097.5714
$ws.Range("K113").Value = 1097.5714
$ws.Range("M113").Value = 1072.4286

# Row 115: At Your Neck and Call / Atrociraptorskin Necklace of Aiming
$ws.Range("H115").Value = 46434
$ws.Range("J115").Value = 46434
$ws.Range("L115").Value = 46434
$ws.Range("N115").Value = -48784

# Row 126: Battered Books / Saiga Leather
$ws.Range("H126").Value = 12862.818
$ws.Range("I126").Value = 12862.818
$ws.Range("K126").Value = 38588.454
$ws.Range("M126").Value = -36118.454

$ws = $wb.Worksheets.Item("WVR")
# Row 100: Of Great Import / Kudzu Thread
$ws.Range("H100").Value = 1769.4
$ws.Range("I100").Value = 1769.4
$ws.Range("K100").Value = 3538.8
$ws.Range("M100").Value = -2997.8

# Row 126: A Polished Purchase / Snow Linen
$ws.Range("H126").Value = 8276.412
$ws.Range("I126").Value = 9899.923000000001
$ws.Range("K126").Value = 29699.769
$ws.Range("M126").Value = -27229.769
